$wb = $excel.ActiveWorkbook

# Add the new "Cases" worksheet after the last existing sheet (Deals)
$lastSheet = $wb.Worksheets($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Cases"

# Header row (yellow-filled, same style as the other sheets' header rows)
$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "status"
$ws.Range("C1").Value = "identifier"
$ws.Range("D1").Value = "type"
$ws.Range("E1").Value = "priority"
$ws.Range("F1").Value = "contact"
$ws.Range("A1:F1").Interior.Color = 65535

# Data rows, entered column by column
$ws.Range("A2").Value = "CaseTitle1"
$ws.Range("A3").Value = "CaseTitle2"

$ws.Range("B2").Value = "Awaiting input"
$ws.Range("B3").Value = "Enquiring"

$ws.Range("C2").Value = "aaaa"
$ws.Range("C3").Value = "bbbb"

$ws.Range("D2").Value = "Business Support"
$ws.Range("D3").Value = "Complaint"

$ws.Range("E2").Value = "High"
$ws.Range("E3").Value = "Low"

$ws.Range("F2").Value = "aaaaa"
$ws.Range("F3").Value = "zzzx"

$ws.Columns("A:F").AutoFit()

# Restore selections/active states to match the final layout
$wb.Worksheets("Companies").Range("E26").Select()
$wb.Worksheets("Deals").Range("A1:XFD1").Select()

$ws.Activate()
$ws.Range("D5").Select()
